$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4178.885805142692
$ws.Range("C3").Value = 4178.885805142692
$ws.Range("C4").Value = 4178.885805142692
$ws.Range("C5").Value = 4142.942492509374
$ws.Range("C6").Value = 4098.780631332637
$ws.Range("C7").Value = 4083.34828679113
$ws.Range("C8").Value = 4061.944200590551
$ws.Range("C9").Value = 4061.944200590551
$ws.Range("C10").Value = 4061.944200590551
$ws.Range("C11").Value = 4061.944200590551
$ws.Range("C12").Value = 4061.944200590551
